$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)
foreach ($r in $rows) {
    $cell = $ws.Range("A$r")
    # Assigning the date-like string directly would make Excel auto-convert
    # it to a date serial number. Force text interpretation by temporarily
    # setting the number format to Text, then restore the original
    # (style-less) formatting so the cell XML matches the source exactly.
    $cell.NumberFormat = "@"
    $cell.Value = "2025/11/10"
    $cell.Style = "Normal"
}
